# Updates cryptocurrency price (column D) and volume-change (column E) values
# on Sheet1, matching the refreshed data pulled by the GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row number, new Price (column D) value, new Volume(1h) (column E) value.
# A value of $null means that column is left unchanged for that row.
$updates = @(
    ,@(2, "62.071.39", $null)
    ,@(3, "3.421.99", "  -0.50%  ")
    ,@(4, $null, "  -0.01%  ")
    ,@(5, $null, "  +0.48%  ")
    ,@(6, "154.12", "  +5.31%  ")
    ,@(7, $null, "  +0.06%  ")
    ,@(8, $null, "  +1.66%  ")
    ,@(9, "8.10", "  +4.64%  ")
    ,@(10, $null, "  +1.27%  ")
    ,@(11, "0.418", "  +3.75%  ")
    ,@(12, "4.008.55", "  -0.43%  ")
    ,@(13, $null, "  +1.04%  ")
    ,@(14, "28.62", "  -0.88%  ")
    ,@(15, $null, "  +0.95%  ")
    ,@(16, "3.416.80", "  -0.27%  ")
    ,@(17, "62.085.47", "  -1.26%  ")
    ,@(18, "6.53", "  +3.03%  ")
    ,@(19, "14.43", "  +0.68%  ")
    ,@(20, "9.00", "  -1.86%  ")
    ,@(21, "383.64", "  -0.28%  ")
    ,@(22, "0.572", "  +2.39%  ")
    ,@(23, "75.90", "  +2.12%  ")
    ,@(24, $null, "  +0.06%  ")
    ,@(25, "3.563.58", "  -0.47%  ")
    ,@(26, $null, "  -0.87%  ")
    ,@(27, $null, "  -1.17%  ")
    ,@(28, $null, "  +1.72%  ")
    ,@(29, $null, "  -0.12%  ")
    ,@(30, $null, "  +1.03%  ")
    ,@(31, "7.87", "  -2.42%  ")
    ,@(32, $null, "  +0.04%  ")
    ,@(33, "23.28", "  +0.17%  ")
    ,@(34, $null, "  +1.73%  ")
    ,@(35, "5.52", "  +4.83%  ")
    ,@(36, $null, "  +0.97%  ")
    ,@(37, $null, "  -1.31%  ")
    ,@(38, "168.74", $null)
    ,@(39, "30.94", "  -2.42%  ")
    ,@(40, "3.458.27", "  -0.52%  ")
    ,@(41, "0.0781", "  +1.75%  ")
    ,@(42, "42.79", "  +1.15%  ")
    ,@(43, "0.781", "  -0.79%  ")
    ,@(44, $null, "  +1.70%  ")
    ,@(45, $null, "  -1.96%  ")
    ,@(46, $null, "  -1.78%  ")
    ,@(47, "2.556.19", "  -0.16%  ")
    ,@(48, "23.33", "  +3.31%  ")
    ,@(49, $null, "  -0.05%  ")
    ,@(50, "2.22", "  -2.19%  ")
    ,@(51, $null, "  +0.01%  ")
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $newPrice = $entry[1]
    $newVolume = $entry[2]

    if ($null -ne $newPrice) {
        # Prefix with an apostrophe so Excel keeps purely-numeric-looking
        # prices (e.g. "154.12") stored as text, matching the sheet's
        # existing inline-string convention for the Price column. Re-apply
        # the Normal style afterwards so the text-forcing quote prefix
        # doesn't leave behind an extra "@ number format" on the cell.
        $ws.Cells.Item($row, 4).Value = "'" + $newPrice
        $ws.Cells.Item($row, 4).Style = "Normal"
    }
    if ($null -ne $newVolume) {
        $ws.Cells.Item($row, 5).Value = $newVolume
    }
}
